$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "from date .. to date" header line with template placeholders
$ws.Range("B6").Value = "Từ ngày `${startDate} đến ngày `${endDate} "

# Fill in the previously-blank signature columns on rows 14-15
$ws.Range("E14").Value = "Kế toán trưởng"
$ws.Range("H14").Value = "Giám đốc"
$ws.Range("E15").Value = "(Ký, họ tên)"
$ws.Range("H15").Value = "(Ký, họ tên, đóng dấu)"

# Match the final selection state recorded in the workbook
$ws.Range("H15:I15").Select()
